# V1.00 As built by PCBway, Aug 7, 2024
#
# Updates the PTHComponents BOM sheet:
#  - Add "Status" column (C) annotations ("?" = unknown/needs sourcing,
#    "dk" = DigiKey, "amazon" = Amazon) for several line items.
#  - Swap the J1 header for the taller 20-pin variant.
#  - Replace the RG174 pigtail assembly with an RG316 pigtail (J2).
#  - Widen the Status column so the new values are readable.
#  - Leave the selection on the last-edited cell (C18).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 : J1 -----------------------------------------------------------
$ws.Range("C2").Value = "?"
$ws.Range("D2").Value = "20-pin 2.54mm pitch tall female header"

# --- Row 3 : J2 -------------------------------------------------------------
$ws.Range("C3").Value = "amazon"
$ws.Range("D3").Value = "RG316"
$ws.Range("F3").Value = "Antenna pigtail"

# --- Row 5 : J5 -------------------------------------------------------------
$ws.Range("C5").Value = "?"

# --- Row 6 : TN1 ------------------------------------------------------------
$ws.Range("C6").Value = "?"

# --- Rows 11-14 : C8, C9, C10, C26 (17M band capacitors) --------------------
$ws.Range("C11").Value = "dk"
$ws.Range("C12").Value = "dk"
$ws.Range("C13").Value = "dk"
$ws.Range("C14").Value = "dk"

# --- Rows 18-19 : Q3, Q4 -----------------------------------------------------
$ws.Range("C18").Value = "dk"
$ws.Range("C19").Value = "dk"

# --- Rows 22-25 : C8, C9, C10, C26 (20M band capacitors) --------------------
$ws.Range("C22").Value = "dk"
$ws.Range("C23").Value = "dk"
$ws.Range("C24").Value = "dk"
$ws.Range("C25").Value = "dk"

# --- Widen the Status column so the new text fits ---------------------------
$ws.Columns.Item(3).ColumnWidth = 16.67

# --- Leave the selection where the author left off --------------------------
[void]$ws.Range("C18").Select()
